# Update the "2024" sheet: a new entry ("dispute" / 2024-09-19 22:33:39) was
# logged at the top of the September activity list (columns R = "...Details",
# S = "...Date", currently holding rows 45-165). Inserting it pushes every
# existing entry in that list down by one row (45->46, 46->47, ..., 165->166).
#
# Separately, the lone "Broadband" label that sits by itself at the bottom of
# column A (row 174) is likewise pushed down by one row, to row 175, which is
# what grows the sheet's used range from A1:Y174 to A1:Y175.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

$colR = 18   # column R - "September_Details"
$colS = 19   # column S - "September_Date"

$firstRow = 45
$lastRow  = 165

# Walk bottom-up so each row's old value is read before it gets overwritten
# by the row below shifting into it.
for ($r = $lastRow; $r -ge $firstRow; $r--) {
    $rVal = $ws.Cells.Item($r, $colR).Value2
    $sVal = $ws.Cells.Item($r, $colS).Value2
    $ws.Cells.Item($r + 1, $colR).Value = $rVal
    $ws.Cells.Item($r + 1, $colS).Value = $sVal
}

# The new, most-recent entry takes the now-vacated top row.
$ws.Cells.Item($firstRow, $colR).Value = "dispute"
$ws.Cells.Item($firstRow, $colS).Value = "2024-09-19 22:33:39"

# The trailing "Broadband" label moves from row 174 to row 175.
$ws.Cells.Item(174, 1).ClearContents()
$ws.Cells.Item(175, 1).Value = "Broadband"
